$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("I18").Copy($ws.Range("G18"))
$ws.Range("I18").ClearContents()
$ws.Columns("H:I").Delete()
Write-Output "G18 = $($ws.Range('G18').Text)"
Write-Output "Dimension check via UsedRange:"
Write-Output $ws.UsedRange.Address()
